$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace old "Recurrent Neural Networks" items with new text, drop wrap style on A2 ---
$ws.Range("A2").Value = "Classical / Statistical Models — Moving Averages, Exponential Smoothing, ARIMA, SARIMA, TBATS"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Improve ARIMA existing model"
$ws.Range("C2").Value = "Selection of scope and direction"

# --- Row 3: new Machine Learning row, taller to fit wrapped long URL text ---
$ws.Range("A3").Value = "Machine Learning — Linear Regression, XGBoost, Random Forest, or any ML model with reduction methods"
$ws.Range("B3").Value = "Dada studying forecasting methods (https://www.datacamp.com/tutorial/tutorial-time-series-forecasting?irclickid=Qu-WEL35QxyIRzmX30wL5WzCUkD2--zut0y1wg0&irgwc=1&utm_medium=affiliate&utm_source=impact&utm_campaign=1310690#what-is-time-series-forecasting-)"
$ws.Range("C3").Value = "Univariate time series forecast of Aquifer Patrignano"
$ws.Rows.Item(3).RowHeight = 29

# --- Row 4: Deep Learning row (keeps its previous 18pt custom height) ---
$ws.Range("A4").Value = "Deep Learning — RNN, LSTM, Autoregressive CNN for univariate time series / multivariate"
$ws.Range("C4").Value = "start ARIMA"

# --- Row 5: new To-Do item ---
$ws.Range("A5").Value = "Products demand (https://www.datacamp.com/courses/forecasting-product-demand-in-r)"

# --- Rows 6-7: new empty wrap-styled placeholder rows ---
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true

# --- Column A is much wider now to hold the new longer To-Do text ---
$ws.Columns.Item(1).ColumnWidth = 75.6

# --- Final selection left on B6 ---
$ws.Range("B6").Select() | Out-Null
